# UserSchedules.xlsx - rework sheets:
#   - drop the empty "Kat" sheet that used to sit first
#   - re-add "Kat" at the end of the tab strip, populated with two rows of
#     class-schedule data
#   - leave "Test2" and "testing" as-is (just shifted earlier in tab order)
#   - make "Test2" (now the first tab) the active tab

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the old, empty "Kat" worksheet.
$wb.Worksheets.Item("Kat").Delete() | Out-Null

# Re-create "Kat" as a brand new sheet after the last existing tab
# (order ends up: Test2, testing, Kat).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$kat = $wb.Worksheets.Add($null, $lastSheet)
$kat.Name = "Kat"

# Row 1
$kat.Range("A1").Value = "HUMA 202  A"
$kat.Range("B1").Value = "CIV/LITERATURE"
$kat.Range("C1").Value = "CIV/LITERATURE"
$kat.Range("D1").Value = "11:00:00"
$kat.Range("E1").Value = "11:50:00"
$kat.Range("F1").Value = "MWF"
$kat.Range("G1").Value = "HAL"

# "302" looks numeric, so force text formatting while writing it, then
# clear the formatting again so the cell keeps using the default style
# (matches the source data being a shared string, not a number).
$kat.Range("H1").NumberFormat = "@"
$kat.Range("H1").Value = "302"
$kat.Range("H1").ClearFormats()

# Row 2
$kat.Range("A2").Value = "HUMA 302  A"
$kat.Range("B2").Value = "MODERN CIV/INTL"
$kat.Range("C2").Value = "MODERN CIV/INTL PERSP"
$kat.Range("D2").Value = "12:00:00"
$kat.Range("E2").Value = "12:50:00"
$kat.Range("F2").Value = "MWF"
$kat.Range("G2").Value = "HAL"

$kat.Range("H2").NumberFormat = "@"
$kat.Range("H2").Value = "304"
$kat.Range("H2").ClearFormats()

# "Test2" is now the first tab - make it the active one.
$wb.Worksheets.Item("Test2").Activate()
